# Update FTSE 100 ticker list: remove "HL / Hargreaves Lansdown / Financial Services"
# row, shifting the subsequent rows (HIK..INF) up by one, and append the new
# "ICP / Intermediate Capital Group / Financial Services" row at the end of the block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New contents for rows 43-50 (Ticker, Company, Sector)
$rows = @(
    @("HIK", "Hikma Pharmaceuticals", "Pharmaceuticals & Biotechnology"),
    @("HWDN", "Howdens Joinery", "Homebuilding & Construction Supplies"),
    @("HSBA", "HSBC", "Banks"),
    @("IHG", "IHG Hotels & Resorts", "Travel & Leisure"),
    @("IMI", "IMI", "Machinery, Tools, Heavy Vehicles, Trains & Ships"),
    @("IMB", "Imperial Brands", "Tobacco"),
    @("INF", "Informa", "Media"),
    @("ICP", "Intermediate Capital Group", "Financial Services")
)

$startRow = 43
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}
